$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two mailto: hyperlinks that used to live on B2 / B4
[void]$ws.Hyperlinks.Delete()

# The refreshed provider list carries new formatting on the e-mail column,
# so start from a clean slate there (keeps B4's old look untouched below)
$ws.Range("B2:B3").ClearFormats()

# Updated provider names
$ws.Range("A2").Value = "ACCESORIOS Y HERRAJES JM SAS"
$ws.Range("A3").Value = "ACDC ELECTRIC SAS"
$ws.Range("A4").ClearContents()

# The e-mail now lines up with row 2 and row 3; row 4 goes blank
$ws.Range("B3").Value = "japsequiposelectricos@gmail.com"
$ws.Range("B4").ClearContents()

# New light-grey Roboto look for the e-mail column (no more hyperlink style)
$ws.Range("B2:B3").Font.Color = 14935011
$ws.Range("B2:B3").Font.Name = "Roboto"

# Stray formatted-but-empty cell left over from the paste
$ws.Range("D2").Value = ""
$ws.Range("D2").Interior.ColorIndex = -4142

# Columns widened to fit the new, longer provider / e-mail text
$ws.Columns("B").AutoFit()
$ws.Columns("C").ColumnWidth = 34.1666667

# Selection ends on B3
[void]$ws.Range("B3").Select()
